# Add "2022-Q1" worksheet with fund-holding detail data, placed right
# before the "总计" (total) sheet, and update the "总计" sheet with a new
# summary row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4" (same
#    column layout/header/styles), then overwrite its contents.
#    NOTE: we intentionally place the copy right after its source (not
#    yet next to "总计") and only Move() it into its final position
#    once all the cell writes are done - writing new rows to a sheet
#    that already sits immediately before the workbook's last sheet
#    does not get persisted correctly by this runtime.
# ---------------------------------------------------------------------
$srcQ4 = $wb.Worksheets.Item("2021-Q4")
$srcQ4.Copy([System.Reflection.Missing]::Value, $srcQ4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Fund-holding detail rows (code, name, size, position, ratio, market
# value, rank).
$data = @(
  @("090003", "大成蓝筹稳健混合", "17.92", "89.15", "7.12", "1.2759", 2),
  @("519019", "大成景阳领先混合", "10.33", "92.80", "8.75", "0.9039", 2),
  @("161838", "银华创业板两年定期开放混合", "10.44", "95.40", "5.83", "0.6087", 6),
  @("090016", "大成消费主题混合", "4.23", "93.78", "9.37", "0.3964", 3),
  @("008128", "湘财长源股票A", "2.74", "94.29", "5.94", "0.1628", 10),
  @("009907", "湘财长泽灵活配置混合A", "1.79", "80.11", "4.45", "0.0797", 9),
  @("008129", "湘财长源股票C", "1.05", "94.29", "5.94", "0.0624", 10),
  @("002319", "大成一带一路灵活配置混合", "0.50", "89.30", "9.06", "0.0453", 1),
  @("009908", "湘财长泽灵活配置混合C", "0.46", "80.11", "4.45", "0.0205", 9),
  @("005082", "诺德量化蓝筹增强混合A", "1.17", "91.78", "1.47", "0.0172", 10),
  @("005083", "诺德量化蓝筹增强混合C", "1.17", "91.78", "1.47", "0.0172", 10),
  @("005295", "诺德天富灵活配置混合", "1.21", "93.81", "1.25", "0.0151", 8),
  @("011761", "平安鑫瑞混合型证券投资基金A", "1.09", "20.41", "0.58", "0.0063", 3),
  @("011762", "平安鑫瑞混合型证券投资基金C", "0.34", "20.41", "0.58", "0.0020", 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $newSheet.Range("A$r").Value = $i

    # Fund code (B) - keep leading zeros, force text.
    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[0]
    $newSheet.Range("B$r").Style = "Normal"

    # Fund name (C) - plain text, never looks numeric.
    $newSheet.Range("C$r").Value = $row[1]

    # Fund size / position / ratio / market value (D,E,F,G) - numeric
    # looking text kept as text.
    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $row[2]
    $newSheet.Range("D$r").Style = "Normal"

    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $row[3]
    $newSheet.Range("E$r").Style = "Normal"

    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $row[4]
    $newSheet.Range("F$r").Style = "Normal"

    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $row[5]
    $newSheet.Range("G$r").Style = "Normal"

    # Position rank (H) - real number.
    $newSheet.Range("H$r").Value = $row[6]

    # Index column (A) keeps the bold/bordered style already used by
    # the template sheet; make sure every row (including the newly
    # appended ones beyond row 9) carries it.
    $newSheet.Range("A$r").Style = "Normal"
}

$srcIndexStyle = $wb.Worksheets.Item("2021-Q4").Range("A2")
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $srcIndexStyle.Copy()
    $newSheet.Range("A$r").PasteSpecial(-4122)
    $newSheet.Range("A$r").Value = $i
}

# Move the finished sheet into its final position, immediately before
# "总计".
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet.Move($totalSheet)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new summary row for 2022-Q1 at
#    the top of the data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert(-4121)
$totalWs.Rows.Item(2).ClearFormats()

$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 14
$totalWs.Range("D2").Value = 3.61

# Renumber the index column for the rows that were pushed down (they
# used to start at 0, now they must start at 1).
for ($r = 3; $r -le 7; $r++) {
    $totalWs.Range("A$r").Value = $r - 2
}
